$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-02-05 Monday" "2024-02-06 Tuesday"
Replace-Text "35÷4=" "27÷8="
Replace-Text "41÷2=" "85÷9="
Replace-Text "93÷8=" "48÷5="
Replace-Text "75÷6=" "79÷3="
Replace-Text "41÷3=" "50÷3="
Replace-Text "88÷3=" "47÷8="
Replace-Text "72÷2=" "23÷7="
Replace-Text "85÷4=" "61÷4="
Replace-Text "50÷7=" "42÷2="
Replace-Text "44÷9=" "80÷8="
Replace-Text "27÷7=" "97÷2="
Replace-Text "59÷6=" "12÷2="
Replace-Text "44÷2=" "99÷4="
Replace-Text "39÷6=" "34÷8="
Replace-Text "96÷8=" "37÷6="
Replace-Text "49÷7=" "52÷7="
Replace-Text "79÷4=" "73÷7="
Replace-Text "55÷8=" "44÷5="
Replace-Text "29÷8=" "66÷7="
Replace-Text "61÷3=" "74÷3="
Replace-Text "67÷7=" "75÷6="
Replace-Text "13÷9=" "72÷2="
Replace-Text "22÷7=" "31÷4="
Replace-Text "46÷2=" "86÷7="
Replace-Text "72÷4=" "26÷3="
